$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pivot Table")
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$s = $chart.SeriesCollection().Item(2)
$s.MarkerStyle = -4142
